$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '71.210.00'
$ws.Range('E2').Value = '  +0.88%  '
$ws.Range('D3').Value = '3.855.10'
$ws.Range('E3').Value = '  +1.08%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = "'697.20"
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +2.37%  '
$ws.Range('D6').Value = "'173.46"
$ws.Range('D6').ClearFormats()
$ws.Range('D7').Value = '3.852.92'
$ws.Range('E7').Value = '  +1.07%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('E9').Value = '  +0.10%  '
$ws.Range('E10').Value = '  +1.16%  '
$ws.Range('D11').Value = "'7.24"
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -0.02%  '
$ws.Range('E12').Value = '  -0.15%  '
$ws.Range('D13').Value = "'0.0000258"
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +5.11%  '
$ws.Range('D14').Value = "'36.43"
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +1.34%  '
$ws.Range('D15').Value = '4.506.17'
$ws.Range('E15').Value = '  +1.12%  '
$ws.Range('D16').Value = '3.855.03'
$ws.Range('E16').Value = '  +1.07%  '
$ws.Range('D17').Value = '71.255.06'
$ws.Range('E17').Value = '  +0.77%  '
$ws.Range('E18').Value = '  +0.07%  '
$ws.Range('D19').Value = "'7.23"
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +0.68%  '
$ws.Range('E20').Value = '  -0.05%  '
$ws.Range('D21').Value = "'11.12"
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -1.08%  '
$ws.Range('D22').Value = "'493.38"
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +3.45%  '
$ws.Range('D23').Value = "'0.723"
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +1.16%  '
$ws.Range('D24').Value = "'85.11"
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +1.97%  '
$ws.Range('E25').Value = '  +1.89%  '
$ws.Range('D26').Value = "'12.29"
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.13%  '
$ws.Range('D27').Value = "'10.62"
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +2.68%  '
$ws.Range('E28').Value = '  +1.66%  '
$ws.Range('D29').Value = '4.009.19'
$ws.Range('E29').Value = '  +1.12%  '
$ws.Range('D30').Value = "'3.19"
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +9.40%  '
$ws.Range('E31').Value = '  -0.01%  '
$ws.Range('D32').Value = "'7.64"
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +3.11%  '
$ws.Range('D33').Value = "'2.29"
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -0.24%  '
$ws.Range('D34').Value = "'29.70"
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +0.43%  '
$ws.Range('E35').Value = '  -1.09%  '
$ws.Range('D36').Value = "'9.28"
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +1.71%  '
$ws.Range('D37').Value = '3.807.18'
$ws.Range('E37').Value = '  +0.99%  '
$ws.Range('E38').Value = '  +0.00%  '
$ws.Range('E39').Value = '  +1.81%  '
$ws.Range('E40').Value = '  +12.52%  '
$ws.Range('D41').Value = "'6.06"
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +1.81%  '
$ws.Range('E42').Value = '  +0.10%  '
$ws.Range('E43').Value = '  +6.66%  '
$ws.Range('E44').Value = '  +0.01%  '
$ws.Range('D46').Value = "'163.98"
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +2.53%  '
$ws.Range('D47').Value = "'0.000308"
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +3.50%  '
$ws.Range('D48').Value = "'48.65"
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +0.93%  '
$ws.Range('D49').Value = "'44.34"
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -3.86%  '
$ws.Range('B50').Value = 'Bittensor'
$ws.Range('C50').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D50').Value = "'419.21"
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +5.05%  '
$ws.Range('B51').Value = 'TheGraph'
$ws.Range('C51').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D51').Value = "'0.303"
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +1.03%  '
